# implemented Axe in screens
# Adds a new "checkAccessibility" / "TC_PM_COA_SEC_ListView_D3" step as row 3
# of the TestSteps sheet (shifting the existing steps down by one row),
# appends a new blank trailing row, and makes TestSteps the active sheet
# (TestData is no longer the active/selected tab).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # TestSteps
$ws2 = $wb.Worksheets.Item(2)   # TestData

# --- 1. Append a brand-new trailing row 16 (mirrors row 15's blank D:F
#        formatting), then shift rows 3..15 down to 4..16 (values +
#        formatting), bottom-up so nothing is overwritten before it has
#        been copied. ---
$ws1.Rows.Item(16).Insert()
$src16 = $ws1.Range("D15:F15")
$dst16 = $ws1.Range("D16:F16")
$src16.Copy()
$dst16.PasteSpecial(-4122)   # xlPasteFormats
$src16.Copy()
$dst16.PasteSpecial(-4163)   # xlPasteValues
$ws1.Application.CutCopyMode = $false

for ($r = 15; $r -ge 3; $r--) {
    $src = $ws1.Range("A" + $r + ":F" + $r)
    $dst = $ws1.Range("A" + ($r + 1) + ":F" + ($r + 1))
    $src.Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats (handles rows that gain newly-used columns)
    $src.Copy()
    $dst.PasteSpecial(-4163)   # xlPasteValues
}
$ws1.Application.CutCopyMode = $false

# --- 2. Populate the new row 3 with the accessibility-check step. ---
$ws1.Range("A3").Value = "checkAccessibility"
$ws1.Range("C3").ClearContents()

# Give B3 the same base border/fill as the rest of the data rows' no-fill
# columns, then layer on the wrap-text / top-left alignment used for the
# long description value.
$ws1.Range("F3").Copy()
$ws1.Range("B3").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Application.CutCopyMode = $false

$ws1.Range("B3").Value = "TC_PM_COA_SEC_ListView_D3"
$ws1.Range("B3").WrapText = $true
$ws1.Range("B3").HorizontalAlignment = -4131  # xlLeft
$ws1.Range("B3").VerticalAlignment = -4160    # xlTop

# --- 3. Make TestSteps the active sheet / selection, and move the
#        TestData selection off of "tab selected". ---
$ws1.Activate()
$ws1.Range("B4").Select()

$ws2.Activate()
$ws2.Range("H8").Select()

$ws1.Activate()
